# Daily User Impact Status - add the 10/21/2025 row of data and move the
# active-cell selection, matching the author's re-upload of the tracker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 was a blank placeholder row (only a date-formatted A11); fill it in
# with the day's totals.
$ws.Range("A11").Value = 45951
$ws.Range("B11").Value = 5598
$ws.Range("C11").Value = 4348
$ws.Range("D11").Value = 3992
$ws.Range("E11").Value = 275
$ws.Range("F11").Value = 42
$ws.Range("G11").Value = 33
$ws.Range("H11").Value = 4
$ws.Range("I11").Value = 2

# The workbook was saved with the active cell on C18 rather than B12.
$ws.Range("C18").Select() | Out-Null
